# "Updated Charges to add DOC - persistent."
#
# The charges table on Sheet1 is kept sorted alphabetically by column A
# ("Name of Offense"). The new charge "Disorderly Conduct - Persistent"
# sorts immediately after the existing "Disorderly Conduct" row (row 7),
# so it is inserted as a new row 8, pushing every following row down by
# one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 8 (shifts rows 8..38 down to 9..39, preserving
# their values/styles/formatting).
$ws.Rows(8).Insert()

# Populate the new row: Name of Offense | Statute | Degree | Offense Type
$ws.Range("A8").Value = "Disorderly Conduct - Persistent"
$ws.Range("B8").Value = "2917.11(A)(1)"
$ws.Range("C8").Value = "M4"
$ws.Range("D8").Value = "Criminal"

# Match the selection left behind in the saved workbook.
$ws.Range("E8").Select() | Out-Null
